# "udpated tutoring to remove darren"
#
# Semantic changes applied:
#  1. Remove Darren Webb's row from the Evening Schedule table (row 21),
#     which shifts Corey Wolf (and the trailing blank spacer rows) up by
#     one row.
#  2. Update the worksheet title to note the new start date.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the title banner text (merged cell A2:L2).
$ws.Range("A2").Value() = "CS320: Software Engineering Tutoring Schedule (tentative - starts Weds, 1-24-18)"

# Delete the "Darren Webb" row from the evening-schedule roster; everything
# below it (Corey Wolf's row and the blank spacer rows) shifts up by one.
$ws.Rows.Item(21).Delete()

# Leave the cursor roughly where the author's saved file had it.
$ws.Range("J17").Select()
